$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: A12's date format changes from date-only to full datetime ---
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 13: new record appended at the bottom of the log ---
$ws.Range("A13").Value = 45856
$ws.Range("A13").NumberFormat = "YYYY-MM-DD"

$ws.Range("B13").Value = "diegonovo"
$ws.Range("C13").Value = "diegonovo"
$ws.Range("D13").Value = "2025-07-18 14:30:11"
$ws.Range("E13").Value = "2025-07-18 14:30:12"
$ws.Range("F13").Value = "2025-07-18 14:30:13"
$ws.Range("G13").Value = "2025-07-18 14:30:14"
$ws.Range("H13").Value = "2025-07-18 14:30:17"
$ws.Range("I13").Value = "2025-07-18 14:30:17"
$ws.Range("J13").Value = "2025-07-18 14:30:23"
$ws.Range("K13").Value = "0:00:01"
$ws.Range("L13").Value = "0:00:01"
$ws.Range("M13").Value = "0:00:12"
$ws.Range("O13").Value = "2025-07-18 14:30:26"
$ws.Range("P13").Value = "2025-07-18 14:30:26"
$ws.Range("Q13").Value = "2025-07-18 14:30:28"
$ws.Range("R13").Value = "2025-07-18 14:30:32"
$ws.Range("S13").Value = "0:00:02"
$ws.Range("T13").Value = "0:00:02"
$ws.Range("U13").Value = "0:00:08"
$ws.Range("V13").Value = "0:00:01"
$ws.Range("W13").Value = "2025-07-18 14:30:24"

Write-Host "Row 13 added"
